$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; existing rows 13-21 shift down to 14-22
$ws.Rows.Item(13).Insert()

# The new row 13 is a new weekly price record for the same market/product,
# carrying the same static attributes as the surrounding rows.
$ws.Cells.Item(13, 1).Value = 10
$ws.Cells.Item(13, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(13, 3).Value = "La Araucanía"
$ws.Cells.Item(13, 4).Value = 44830
$ws.Cells.Item(13, 5).Value = 9
$ws.Cells.Item(13, 6).Value = "Fruta"
$ws.Cells.Item(13, 7).Value = 100108
$ws.Cells.Item(13, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(13, 9).Value = 100108001
$ws.Cells.Item(13, 10).Value = "Guayaba"
$ws.Cells.Item(13, 11).Value = "Sin especificar"
$ws.Cells.Item(13, 12).Value = "Primera"
$ws.Cells.Item(13, 13).Value = 50
$ws.Cells.Item(13, 14).Value = 2500
$ws.Cells.Item(13, 15).Value = 2500
$ws.Cells.Item(13, 16).Value = 2500
$ws.Cells.Item(13, 17).Value = "$/kilo"
$ws.Cells.Item(13, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(13, 19).Value = 2500
$ws.Cells.Item(13, 20).Value = 1
